$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @(44319, 44320, 44321)
$bvals = @(1, 0, 0)
$cvals = @(3, 3, 3)
$dvals = @(112.4016485575122, 112.4016485575122, 112.4016485575122)

for ($i = 0; $i -lt 3; $i++) {
    $row = 245 + $i

    # Copy formatting from the last existing row (244) for column A (date style)
    $ws.Cells.Item(244, 1).Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4122)
    $ws.Cells.Item($row, 1).Value = $dates[$i]

    $ws.Cells.Item($row, 2).Value = $bvals[$i]
    $ws.Cells.Item($row, 3).Value = $cvals[$i]
    $ws.Cells.Item($row, 4).Value = $dvals[$i]
}

$excel.CutCopyMode = 0
